# The table cell that used to read "productId" (stored as the two runs
# "productI" + "d") should read "productId" -> "productionId", i.e. the
# three letters "ion" get inserted between "product" and "Id".
$d = $word.ActiveDocument

$found = $d.Content
$ok = $found.Find.Execute("productId", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)

if ($ok) {
    $wordStart = $found.Start

    # "product" is the first 7 characters of the match; the insertion point
    # sits right after it (and right before the "I" of the old "productI").
    $splitPoint = $wordStart + 7

    $insPt = $d.Range($splitPoint, $splitPoint)
    $insPt.InsertBefore("ion")

    # Nudge the freshly inserted "ion" with a formatting round-trip so it is
    # kept as its own run instead of being silently re-absorbed into the
    # neighbouring text.
    $rIon = $d.Range($splitPoint, $splitPoint + 3)
    $rIon.Bold = 1
    $rIon.Bold = 0

    # Same treatment for the trailing "I" so it stays separate from "d".
    $rTailI = $d.Range($splitPoint + 3, $splitPoint + 4)
    $rTailI.Bold = 1
    $rTailI.Bold = 0
}
